$p = $ppt.ActivePresentation

# --- Slide 4: merge the two trailing runs of the "Fixed point ..." bullet
# into a single run (keeps the leading "F" run untouched). ---
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
$para4 = $sh4.TextFrame.TextRange.Paragraphs(4)
$run2 = $para4.Runs(2)
$run3 = $para4.Runs(3)
$run2.Text = "ixed point representations throughout system"
$run3.Text = ""

# --- Slide 5: reposition/resize the body placeholder (PDF version add). ---
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(2)
$sh5.Left = 18
$sh5.Top = 42
$sh5.Width = 719.9716535433071
$sh5.Height = 455.98110236220475
